$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing row 102 values (B..G) with refreshed data
$ws.Range("B102").Value = 54616
$ws.Range("C102").Value = 10810
$ws.Range("D102").Value = 43806
$ws.Range("E102").Value = 17252
$ws.Range("F102").Value = 15434
$ws.Range("G102").Value = 56433

# Add new row 103 with the new quarter entry.
# A103 must hold the literal text "01-04-2021" (not be auto-converted to a
# date serial). Writing it via Formula as a quoted string literal and then
# converting the formula to a static value via copy/paste-special keeps it
# as plain text without altering the cell's style.
$ws.Range("A103").Formula = "=""01-04-2021"""
$ws.Range("A103").Copy()
$ws.Range("A103").PasteSpecial(-4163)
$excel.CutCopyMode = $false

$ws.Range("B103").Value = 55836
$ws.Range("C103").Value = 11578
$ws.Range("D103").Value = 44258
$ws.Range("E103").Value = 17716
$ws.Range("F103").Value = 16126
$ws.Range("G103").Value = 57427
